# Update "想去人数" (F column) counts on both the "展览" sheet and the
# "全部类型" sheet (which mirrors the same rows) to reflect newly
# generated output data.

$wb = $excel.ActiveWorkbook

# Row (in each sheet) -> new F-column value
$updates = @{
    2  = 195
    4  = 12654
    5  = 1285
    6  = 156
    10 = 210
    11 = 457
    17 = 5438
    22 = 127
    23 = 82
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
